$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row at 95, pushing the old row 95 (the footnote row) down to 96,
# and inheriting the same formatting as the row above it (row 94).
$ws.Rows.Item(95).Insert()

# Fill in the new daily-consultation data row (date 2020-04-29 and its counts).
$ws.Range("A95").Value = 43950
$ws.Range("B95").Value = 396
$ws.Range("C95").Value = 31510
$ws.Range("D95").Value = 0
$ws.Range("E95").Value = 6664

# Update the selection to match the new state.
$null = $ws.Range("A95").Select()

# Update the print area to cover the new last row (E98, since the footnote moved to row 96
# and the page's blank trailer rows now extend one row further).
$ws.PageSetup.PrintArea = "相談件数!`$A`$1:`$E`$98"
